$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text values in the Price column are not auto-converted to numbers
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Update Price (column D) values
$ws.Range("D2").Value = "59.115.68"
$ws.Range("D3").Value = "2.592.41"
$ws.Range("D5").Value = "520.17"
$ws.Range("D6").Value = "139.10"
$ws.Range("D7").Value = "0.997"
$ws.Range("D8").Value = "0.566"
$ws.Range("D9").Value = "2.612.23"
$ws.Range("D12").Value = "0.331"
$ws.Range("D14").Value = "3.054.81"
$ws.Range("D15").Value = "59.115.01"
$ws.Range("D16").Value = "20.35"
$ws.Range("D17").Value = "2.588.66"
$ws.Range("D18").Value = "0.0000132"
$ws.Range("D19").Value = "338.79"
$ws.Range("D20").Value = "4.31"
$ws.Range("D21").Value = "10.16"
$ws.Range("D22").Value = "6.48"
$ws.Range("D23").Value = "0.997"
$ws.Range("D26").Value = "0.403"
$ws.Range("D28").Value = "7.03"
$ws.Range("D29").Value = "0.998"
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("D32").Value = "18.79"
$ws.Range("D35").Value = "3.98"
$ws.Range("D36").Value = "1.12"
$ws.Range("D37").Value = "36.29"
$ws.Range("D38").Value = "1.45"
$ws.Range("D39").Value = "0.826"
$ws.Range("D40").Value = "0.826"
$ws.Range("D43").Value = "274.85"
$ws.Range("D45").Value = "0.590"
$ws.Range("D46").Value = "0.0949"
$ws.Range("D47").Value = "0.0520"
$ws.Range("D48").Value = "18.49"
$ws.Range("D49").Value = "1.983.69"
$ws.Range("D50").Value = "4.59"

# Restore default style on the Price column so no formatting side effects remain
$priceRange.Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +6.00%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("E31").Value = "  -3.73%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("E43").Value = "  +5.88%  "
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("E51").Value = "  -0.61%  "

Write-Host "Updated cryptos list values"
